$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F/G columns for rows 5-14 (block A, User "A") and rows 20-29 (block B, User "B")
# These cells hold test-case labels and JSON payloads that were restructured
# from per-item inputs (CreateInput0a/0b, CreateInput1a/1b, UpdateInput0/1) into
# array-based inputs (CreateInputs0/1, UpdateInputs0/1) plus new DeleteIds0/1 rows.

# Block A
$ws.Range("F5").Value = 'CreateInputs0'
$ws.Range("G5").Value = '[{"Id":-999005,"Name":"Marcia"},{"Id":-999007,"Name":"Bobby"}]'
$ws.Range("F6").Value = 'CreateInputs1'
$ws.Range("G6").Value = '[{"Id":-999005,"Name":"Peter"},{"Id":-999006,"Name":"Jan"}]'
$ws.Range("F7").Value = 'CreateExpected0'
$ws.Range("G7").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Marcia", "SysUser":"moe"}, {"Id":-999007, "Name":"Bobby", "SysUser":"moe"}]'
$ws.Range("F8").Value = 'CreateExpected1'
$ws.Range("G8").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Peter", "SysUser":"larry"}, {"Id":-999006, "Name":"Jan", "SysUser":"larry"}]'
$ws.Range("F9").Value = 'UpdateInputs0'
$ws.Range("G9").Value = '[{"Id":-999005,"Name":"Alice"}]'
$ws.Range("F10").Value = 'UpdateInputs1'
$ws.Range("G10").Value = '[{"Id":-999006,"Name":"Cindy"}]'
$ws.Range("F11").Value = 'UpdateExpected0'
$ws.Range("G11").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Alice", "SysUser":"moe"}, {"Id":-999007, "Name":"Bobby", "SysUser":"moe"}]'
$ws.Range("F12").Value = 'UpdateExpected1'
$ws.Range("G12").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Peter", "SysUser":"larry"}, {"Id":-999006, "Name":"Cindy", "SysUser":"larry"}]'
$ws.Range("F13").Value = 'DeleteIds0'
$ws.Range("G13").Value = '[-999005]'
$ws.Range("F14").Value = 'DeleteIds1'
$ws.Range("G14").Value = '[-999005]'

# Block B
$ws.Range("F20").Value = 'CreateInput0'
$ws.Range("G20").Value = '[{"Id":-999005,"Name":"Marcia"},{"Id":-999007,"Name":"Bobby"}]'
$ws.Range("F21").Value = 'CreateInput1'
$ws.Range("G21").Value = '[{"Id":-999005,"Name":"Peter"},{"Id":-999006,"Name":"Jan"}]'
$ws.Range("F22").Value = 'CreateExpected0'
$ws.Range("G22").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Marcia", "SysUser":"moe"}, {"Id":-999007, "Name":"Bobby", "SysUser":"moe"}]'
$ws.Range("F23").Value = 'CreateExpected1'
$ws.Range("G23").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Peter", "SysUser":"larry"}, {"Id":-999006, "Name":"Jan", "SysUser":"larry"}]'
$ws.Range("F24").Value = 'UpdateInput0'
$ws.Range("G24").Value = '[{"Id":-999005,"Name":"Alice"}]'
$ws.Range("F25").Value = 'UpdateInput1'
$ws.Range("G25").Value = '[{"Id":-999006,"Name":"Cindy"}]'
$ws.Range("F26").Value = 'UpdateExpected0'
$ws.Range("G26").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Alice", "SysUser":"moe"}, {"Id":-999007, "Name":"Bobby", "SysUser":"moe"}]'
$ws.Range("F27").Value = 'UpdateExpected1'
$ws.Range("G27").Value = '[{"Id":-999001, "Name":"Mike", "SysUser":"jack@hill.org"}, {"Id":-999002, "Name":"Carol", "SysUser":"jill@hill.org"}, {"Id":-999003, "Name":"Greg", "SysUser":"jack@hill.org"}, {"Id":-999005, "Name":"Peter", "SysUser":"larry"}, {"Id":-999006, "Name":"Cindy", "SysUser":"larry"}]'
$ws.Range("F28").Value = 'DeleteIds0'
$ws.Range("G28").Value = '[-999005]'
$ws.Range("F29").Value = 'DeleteIds1'
$ws.Range("G29").Value = '[-999005]'

# Update the active selection to match the saved cursor position
$ws.Range("E33").Select()
